# Weekly Fruit/Vegetable price update for "Fruta, Vega Central Mapocho de
# Santiago - Pomelo": a new daily record is inserted at row 33 (pushing the
# existing rows 33-76 down to 34-77), and the worksheet's used range grows
# from A1:T76 to A1:T77 accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the current row 33; this shifts every row
# from 33 downward by one position (old row 33 becomes row 34, ..., old
# row 76 becomes row 77).
$ws.Rows.Item(33).Insert()

# Populate the newly inserted row 33 with the new record.
$ws.Range("A33").Value = 9
$ws.Range("B33").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C33").Value = "Metropolitana"
$ws.Range("D33").Value = 44771
$ws.Range("E33").Value = 13
$ws.Range("F33").Value = "Fruta"
$ws.Range("G33").Value = 100102
$ws.Range("H33").Value = "Cítricos"
$ws.Range("I33").Value = 100102006
$ws.Range("J33").Value = "Pomelo"
$ws.Range("K33").Value = "Start Ruby"
$ws.Range("L33").Value = "Primera"
$ws.Range("M33").Value = 250
$ws.Range("N33").Value = 8500
$ws.Range("O33").Value = 8500
$ws.Range("P33").Value = 8500
$ws.Range("Q33").Value = "$/caja 14 kilos"
$ws.Range("R33").Value = "Región Metropolitana"
$ws.Range("S33").Value = 607
$ws.Range("T33").Value = 14
